# This script re-orders the structural-match table on the active worksheet.
# The table (rows 2..97, columns A..H) is made up of 16 "blocks" of 6 rows
# each. Within a block, columns A-D (the "left" entity) are identical across
# all 6 rows, while columns E-H (the "right" entity match) differ. The edit:
#   1) re-orders the 16 blocks (the left entities) into a new sequence, and
#   2) re-orders the 6 right-hand rows inside every block into a new,
#      uniform sequence (keyed by the time-unit suffix of the right_iri).
# This reproduces the axioms added "to allow for matching b/w bfo & ies."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 97
$numRows = $lastRow - $firstRow + 1
$blockSize = 6
$numBlocks = $numRows / $blockSize

# Read the whole data block (A2:H97) as a 1-based [row, col] array.
$data = $ws.Range("A$firstRow`:H$lastRow").Value2

# --- helper: pull the short id / fragment off the end of an iri ---------
function Get-UrlSuffix($url, $sep) {
    $parts = $url.Split($sep)
    return $parts[$parts.Count - 1]
}

# --- group the existing rows into blocks --------------------------------
# Each block remembers: its left-entity key (A col suffix), its A-D values
# (taken from the block's first row), and a lookup (by right_iri suffix)
# of its 6 rows' E-H values.
# NOTE: nested hashtables are mutated through the full "$blocksByKey[$key].Rows[...]"
# path rather than via an intermediate variable, since an intermediate
# variable holds a copy rather than a live reference in this environment.
$blocksByKey = @{}
$blockOrderOriginal = New-Object System.Collections.Generic.List[string]

for ($b = 0; $b -lt $numBlocks; $b++) {
    $r0 = $b * $blockSize + 1   # 1-based row offset into $data for this block's first row

    $leftIri = $data[$r0, 1]
    $key = Get-UrlSuffix $leftIri "/"

    if (-not $blocksByKey.ContainsKey($key)) {
        [void]$blockOrderOriginal.Add($key)

        $blocksByKey[$key] = @{
            A = $data[$r0, 1]
            B = $data[$r0, 2]
            C = $data[$r0, 3]
            D = $data[$r0, 4]
            Rows = @{}
        }
    }

    for ($i = 0; $i -lt $blockSize; $i++) {
        $r = $r0 + $i
        $rightIri = $data[$r, 5]
        $rightKey = Get-UrlSuffix $rightIri "#"
        $rowVals = @($data[$r,1], $data[$r,2], $data[$r,3], $data[$r,4], $data[$r,5], $data[$r,6], $data[$r,7], $data[$r,8])
        $blocksByKey[$key].Rows[$rightKey] = $rowVals
    }
}

# --- the new ordering of the 16 left-entity blocks -----------------------
$newBlockOrder = @(
    "ont00001166",
    "ont00000800",
    "ont00000810",
    "ont00000992",
    "ont00000085",
    "ont00000619",
    "ont00000211",
    "ont00000223",
    "ont00000063",
    "ont00000225",
    "ont00001154",
    "ont00001206",
    "ont00001058",
    "ont00000832",
    "ont00000329",
    "ont00000359"
)

# --- the new ordering of the 6 right-entity rows inside every block -------
$newRowOrder = @("MonthOfYear", "TemporalPosition", "Year", "Duration", "TimePosition", "January")

# --- build the replacement array and write it back ------------------------
$newData = New-Object 'object[,]' $numRows, 8

for ($b = 0; $b -lt $numBlocks; $b++) {
    $key = $newBlockOrder[$b]
    $block = $blocksByKey[$key]

    for ($i = 0; $i -lt $blockSize; $i++) {
        $rightKey = $newRowOrder[$i]
        $rowVals = $block.Rows[$rightKey]

        $destRow = $b * $blockSize + $i   # 0-based row index into $newData
        for ($c = 0; $c -lt 8; $c++) {
            $newData[$destRow, $c] = $rowVals[$c]
        }
    }
}

$ws.Range("A$firstRow`:H$lastRow").Value2 = $newData
